$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the small F2:G2 pair up to F1:G1
$ws.Range("F2:G2").Cut($ws.Range("F1"))

# Move F3 up to F2
$ws.Range("F3").Cut($ws.Range("F2"))

# Move the D6:E10 lookup table to K1:L5 (column range A:B parsing test data)
$ws.Range("D6:E10").Cut($ws.Range("K1"))

# Update the selection/active cell to A11
$ws.Range("A11").Select()
